$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write the new quarter's label as a literal text string (not a date).
# Assigning the text directly to .Value would trigger Excel's smart
# date-recognition for "dd-mm-yyyy"-shaped input, so instead enter it as a
# text formula first and then collapse the formula down to its literal
# text result via copy / paste-special-values - this avoids creating any
# extra number-format/style entries.
$ws.Range("A54").Formula = "=""01-04-2021"""
$ws.Range("A54").Copy()
$ws.Range("A54").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B54").Value = 24.6
$ws.Range("C54").Value = 0.3
$ws.Range("D54").Value = 24.9
